# Update product row 2 (Modelo Prueba sheet): Sizes / Stocks / Codebars columns
# were simplified from multi-variant strings to single-variant strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modelo Prueba")

$ws.Range("G2").Value = "L,M ; L,M"
$ws.Range("H2").Value = "1;1"
$ws.Range("I2").Value = "20310135;12121212"

# Reflect the author's last on-screen selection/scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F14").Select()
